$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

Set-TextValue $ws.Range('D2') '61.578.68'
Set-TextValue $ws.Range('E2') '  +11.70%  '
Set-TextValue $ws.Range('D3') '2.610.91'
Set-TextValue $ws.Range('E3') '  +11.49%  '
Set-TextValue $ws.Range('D4') '0.998'
Set-TextValue $ws.Range('E4') '  -0.19%  '
Set-TextValue $ws.Range('D5') '516.54'
Set-TextValue $ws.Range('E5') '  +8.74%  '
Set-TextValue $ws.Range('D6') '159.56'
Set-TextValue $ws.Range('E6') '  +10.35%  '
Set-TextValue $ws.Range('D7') '0.994'
Set-TextValue $ws.Range('E7') '  -0.45%  '
Set-TextValue $ws.Range('D8') '0.611'
Set-TextValue $ws.Range('E8') '  +0.17%  '
Set-TextValue $ws.Range('D9') '2.650.67'
Set-TextValue $ws.Range('E9') '  +13.07%  '
Set-TextValue $ws.Range('D10') '6.16'
Set-TextValue $ws.Range('E10') '  +13.23%  '
Set-TextValue $ws.Range('D11') '0.106'
Set-TextValue $ws.Range('E11') '  +10.85%  '
Set-TextValue $ws.Range('D12') '0.349'
Set-TextValue $ws.Range('E12') '  +7.38%  '
Set-TextValue $ws.Range('D13') '0.127'
Set-TextValue $ws.Range('E13') '  +1.49%  '
Set-TextValue $ws.Range('D14') '3.061.00'
Set-TextValue $ws.Range('E14') '  +11.41%  '
Set-TextValue $ws.Range('D15') '60.928.79'
Set-TextValue $ws.Range('E15') '  +10.49%  '
Set-TextValue $ws.Range('D16') '22.45'
Set-TextValue $ws.Range('E16') '  +12.63%  '
Set-TextValue $ws.Range('E17') '  +9.82%  '
Set-TextValue $ws.Range('D18') '2.628.86'
Set-TextValue $ws.Range('E18') '  +12.01%  '
Set-TextValue $ws.Range('D19') '4.84'
Set-TextValue $ws.Range('E19') '  +6.61%  '
Set-TextValue $ws.Range('D20') '355.36'
Set-TextValue $ws.Range('E20') '  +13.39%  '
Set-TextValue $ws.Range('D21') '10.62'
Set-TextValue $ws.Range('E21') '  +11.23%  '
Set-TextValue $ws.Range('D22') '6.21'
Set-TextValue $ws.Range('E22') '  +10.31%  '
Set-TextValue $ws.Range('D23') '1.00'
Set-TextValue $ws.Range('E23') '  +0.19%  '
Set-TextValue $ws.Range('D24') '61.11'
Set-TextValue $ws.Range('E24') '  +8.91%  '
Set-TextValue $ws.Range('D25') '0.427'
Set-TextValue $ws.Range('E25') '  +8.70%  '
Set-TextValue $ws.Range('E26') '  +12.33%  '
Set-TextValue $ws.Range('D27') '2.699.64'
Set-TextValue $ws.Range('E27') '  +10.89%  '
Set-TextValue $ws.Range('D28') '0.988'
Set-TextValue $ws.Range('E28') '  -1.11%  '
Set-TextValue $ws.Range('D29') '0.0₃0870'
Set-TextValue $ws.Range('E29') '  +17.80%  '
Set-TextValue $ws.Range('D30') '7.61'
Set-TextValue $ws.Range('E30') '  +8.09%  '
Set-TextValue $ws.Range('D31') '0.997'
Set-TextValue $ws.Range('E31') '  -0.23%  '
Set-TextValue $ws.Range('D32') '19.81'
Set-TextValue $ws.Range('E32') '  +9.55%  '
Set-TextValue $ws.Range('D33') '157.89'
Set-TextValue $ws.Range('E33') '  +8.24%  '
Set-TextValue $ws.Range('E34') '  +8.69%  '
Set-TextValue $ws.Range('D35') '5.62'
Set-TextValue $ws.Range('E35') '  +10.97%  '
Set-TextValue $ws.Range('E36') '  +11.95%  '
Set-TextValue $ws.Range('D37') '4.02'
Set-TextValue $ws.Range('E37') '  +11.85%  '
Set-TextValue $ws.Range('D38') '0.881'
Set-TextValue $ws.Range('E38') '  +9.25%  '
Set-TextValue $ws.Range('D39') '1.50'
Set-TextValue $ws.Range('E39') '  +14.13%  '
Set-TextValue $ws.Range('D40') '308.54'
Set-TextValue $ws.Range('E40') '  +24.09%  '
Set-TextValue $ws.Range('D41') '3.79'
Set-TextValue $ws.Range('E41') '  +12.64%  '
Set-TextValue $ws.Range('D42') '35.67'
Set-TextValue $ws.Range('E42') '  +6.11%  '
Set-TextValue $ws.Range('B43') 'SuiNetwork'
Set-TextValue $ws.Range('C43') 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
Set-TextValue $ws.Range('D43') '0.818'
Set-TextValue $ws.Range('E43') '  +32.46%  '
Set-TextValue $ws.Range('B44') 'Mantle'
Set-TextValue $ws.Range('C44') 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue $ws.Range('D44') '0.640'
Set-TextValue $ws.Range('E44') '  +11.23%  '
Set-TextValue $ws.Range('D45') '0.0580'
Set-TextValue $ws.Range('E45') '  +12.97%  '
Set-TextValue $ws.Range('E46') '  +0.97%  '
Set-TextValue $ws.Range('D47') '20.08'
Set-TextValue $ws.Range('E47') '  +21.09%  '
Set-TextValue $ws.Range('B48') 'FirstDigitalUSD'
Set-TextValue $ws.Range('C48') 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue $ws.Range('D48') '0.991'
Set-TextValue $ws.Range('E48') '  -0.77%  '
Set-TextValue $ws.Range('B49') 'RenderToken'
Set-TextValue $ws.Range('C49') 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws.Range('D49') '5.03'
Set-TextValue $ws.Range('E49') '  +15.89%  '
Set-TextValue $ws.Range('D50') '0.0240'
Set-TextValue $ws.Range('E50') '  +8.81%  '
Set-TextValue $ws.Range('D51') '2.024.99'
Set-TextValue $ws.Range('E51') '  +13.33%  '
